# Generate Report for Handoff
# Updates the localization-status workbook: the previous handoff round
# (1f4f15df...md / eae258fe...md) is replaced by a new handoff round
# (17626ed7...md plus two new dependent .png assets), with refreshed
# status/timestamps, on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet 1: Overview
# -------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Add row 4 (copy formatting from row 3 which already has the right styles)
$ws1.Rows.Item(3).Copy()
$ws1.Rows.Item(4).Insert(-4121)

$ws1.Range("A2").Value2 = "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"
$ws1.Range("B2").Value2 = "Ready for handoff"
$ws1.Range("C2").Value2 = "Ready for handoff"
$ws1.Range("D2").Value2 = "2016-03-20 05:13:24"

$ws1.Range("A3").Value2 = "64ddc91b-10d8-4421-bfb7-d2f141613805.png"
$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"
$ws1.Range("D3").Value2 = "2016-03-20 05:13:24"

$ws1.Range("A4").Value2 = "81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png"
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"
$ws1.Range("D4").Value2 = "2016-03-20 05:13:24"

# Rebuild hyperlinks for column A
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md", [Type]::Missing, [Type]::Missing, "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/64ddc91b-10d8-4421-bfb7-d2f141613805.png", [Type]::Missing, [Type]::Missing, "64ddc91b-10d8-4421-bfb7-d2f141613805.png") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png", [Type]::Missing, [Type]::Missing, "81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png") | Out-Null

# -------------------------------------------------------------------
# Sheet 2: zh-cn
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Add row 4 (copy formatting from row 3)
$ws2.Rows.Item(3).Copy()
$ws2.Rows.Item(4).Insert(-4121)

# Row 2 (handoff file)
$ws2.Range("A2").Value2 = "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"
$ws2.Range("B2").Value2 = ".md"
$ws2.Range("C2").Value2 = "Ready for handoff"
$ws2.Range("D2").Value2 = "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.43651e32d2e57caae6fd88d114a1c1312ef793c4.zh-cn.xlf"
$ws2.Range("E2").Value2 = "2016-03-20 05:13:16"
$ws2.Range("F2").Clear()
$ws2.Range("G2").Clear()
$ws2.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws2.Range("J2").Value2 = "Include"

# Row 3 (first dependent png)
$ws2.Range("A3").Value2 = "64ddc91b-10d8-4421-bfb7-d2f141613805.png"
$ws2.Range("B3").Value2 = ".png"
$ws2.Range("C3").Value2 = "Ready for handoff"
$ws2.Range("D3").Value2 = "8b0378f124537e54d5c91bd0e5fca78f3555d9b7.png"
$ws2.Range("E3").Value2 = "2016-03-20 05:13:16"
$ws2.Range("F3").Clear()
$ws2.Range("G3").Clear()
$ws2.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws2.Range("J3").Value2 = "IsDependency"
$ws2.Range("K3").Value2 = "e2e\17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"

# Row 4 (second dependent png, new row)
$ws2.Range("A4").Value2 = "81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png"
$ws2.Range("B4").Value2 = ".png"
$ws2.Range("C4").Value2 = "Ready for handoff"
$ws2.Range("D4").Value2 = "9651530ea8918632dd986d0731aa2a8a4be71ed2.png"
$ws2.Range("E4").Value2 = "2016-03-20 05:13:16"
$ws2.Range("F4").Clear()
$ws2.Range("G4").Clear()
$ws2.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("J4").Value2 = "IsDependency"
$ws2.Range("K4").Value2 = "e2e\17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"

# Rebuild hyperlinks for columns A and D
$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md", [Type]::Missing, [Type]::Missing, "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a642349d98c04b247ab23b98000acea532a45c67/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.43651e32d2e57caae6fd88d114a1c1312ef793c4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.43651e32d2e57caae6fd88d114a1c1312ef793c4.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/64ddc91b-10d8-4421-bfb7-d2f141613805.png", [Type]::Missing, [Type]::Missing, "64ddc91b-10d8-4421-bfb7-d2f141613805.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a642349d98c04b247ab23b98000acea532a45c67/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/8b0378f124537e54d5c91bd0e5fca78f3555d9b7.png", [Type]::Missing, [Type]::Missing, "8b0378f124537e54d5c91bd0e5fca78f3555d9b7.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png", [Type]::Missing, [Type]::Missing, "81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a642349d98c04b247ab23b98000acea532a45c67/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/9651530ea8918632dd986d0731aa2a8a4be71ed2.png", [Type]::Missing, [Type]::Missing, "9651530ea8918632dd986d0731aa2a8a4be71ed2.png") | Out-Null

# -------------------------------------------------------------------
# Sheet 3: de-de
# -------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Add row 4 (copy formatting from row 3)
$ws3.Rows.Item(3).Copy()
$ws3.Rows.Item(4).Insert(-4121)

# Row 2 (handoff file)
$ws3.Range("A2").Value2 = "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"
$ws3.Range("B2").Value2 = ".md"
$ws3.Range("C2").Value2 = "Ready for handoff"
$ws3.Range("D2").Value2 = "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.43651e32d2e57caae6fd88d114a1c1312ef793c4.de-de.xlf"
$ws3.Range("E2").Value2 = "2016-03-20 05:13:24"
$ws3.Range("F2").Clear()
$ws3.Range("G2").Clear()
$ws3.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws3.Range("J2").Value2 = "Include"

# Row 3 (first dependent png)
$ws3.Range("A3").Value2 = "64ddc91b-10d8-4421-bfb7-d2f141613805.png"
$ws3.Range("B3").Value2 = ".png"
$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("D3").Value2 = "8b0378f124537e54d5c91bd0e5fca78f3555d9b7.png"
$ws3.Range("E3").Value2 = "2016-03-20 05:13:24"
$ws3.Range("F3").Clear()
$ws3.Range("G3").Clear()
$ws3.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws3.Range("J3").Value2 = "IsDependency"
$ws3.Range("K3").Value2 = "e2e\17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"

# Row 4 (second dependent png, new row)
$ws3.Range("A4").Value2 = "81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png"
$ws3.Range("B4").Value2 = ".png"
$ws3.Range("C4").Value2 = "Ready for handoff"
$ws3.Range("D4").Value2 = "9651530ea8918632dd986d0731aa2a8a4be71ed2.png"
$ws3.Range("E4").Value2 = "2016-03-20 05:13:24"
$ws3.Range("F4").Clear()
$ws3.Range("G4").Clear()
$ws3.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("J4").Value2 = "IsDependency"
$ws3.Range("K4").Value2 = "e2e\17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md"

# Rebuild hyperlinks for columns A and D
$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md", [Type]::Missing, [Type]::Missing, "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5887ac1ee6ec63daf730237674f20e8c384763bf/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.43651e32d2e57caae6fd88d114a1c1312ef793c4.de-de.xlf", [Type]::Missing, [Type]::Missing, "17626ed7-f0f4-4ea0-9533-7a2a190ef8ac.43651e32d2e57caae6fd88d114a1c1312ef793c4.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/64ddc91b-10d8-4421-bfb7-d2f141613805.png", [Type]::Missing, [Type]::Missing, "64ddc91b-10d8-4421-bfb7-d2f141613805.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5887ac1ee6ec63daf730237674f20e8c384763bf/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/8b0378f124537e54d5c91bd0e5fca78f3555d9b7.png", [Type]::Missing, [Type]::Missing, "8b0378f124537e54d5c91bd0e5fca78f3555d9b7.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b3d2a2bea39b6bcc71acaa035ecca056766bba60/e2e/81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png", [Type]::Missing, [Type]::Missing, "81d7aded-ebd7-4896-a6c3-c7d276ca09d1.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5887ac1ee6ec63daf730237674f20e8c384763bf/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/9651530ea8918632dd986d0731aa2a8a4be71ed2.png", [Type]::Missing, [Type]::Missing, "9651530ea8918632dd986d0731aa2a8a4be71ed2.png") | Out-Null
